$d = $word.ActiveDocument

# Footer (index 1 = "Approved by: Delivery Manager" footer / footer2.xml):
#   Pearson logo inline picture: name="image1.png" -> name="image2.png"
$d.Sections(1).Footers.Item(1).Range.InlineShapes.Item(1).Name = "image2.png"

# Footer (index 2 = "Authorised by: Head of BTEC Assessment" footer / footer1.xml):
#   Pearson logo inline picture: name="image1.png" -> name="image2.png"
$d.Sections(1).Footers.Item(2).Range.InlineShapes.Item(1).Name = "image2.png"

# Header (index 2 = first-page header / header1.xml):
#   BTEC logo inline picture: name="image2.jpg" -> name="image1.jpg"
$d.Sections(1).Headers.Item(2).Range.InlineShapes.Item(1).Name = "image1.jpg"
